$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6077240705490112
$ws.Range("B1").Value = 1.380396723747253
$ws.Range("C1").Value = 5.642730236053467
$ws.Range("D1").Value = 1.624944090843201
$ws.Range("E1").Value = 1.079240202903748
